# Ajout Timer et modification du zoom en fin de partie
# -> Adds a "Commentaire" column (D) to the Backlog sheet, narrows column C
#    (moving its old width to the new column D), fills in DRI / OK for two
#    tasks that previously had blank Qui / Fait, and resets the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Commentaire" header in D1, formatted like the other header cells ---
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "Commentaire"

# --- Give the rest of column D (D2:D29) the same bordered style used by B/C ---
$ws.Range("B4").Copy() | Out-Null
$ws.Range("D2:D29").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Fill DRI / OK for "Réaliser animation player 2" (row 12) and
#     "Affichage victoire joueur" (row 15), matching the existing pattern ---
$ws.Range("B12").Value = "DRI"
$ws.Range("C12").Value = "OK"
$ws.Range("B15").Value = "DRI"
$ws.Range("C15").Value = "OK"

# --- Resize columns: C becomes narrow, D takes over the previous C width ---
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
$ws.Columns.Item(4).ColumnWidth = 36.5

# --- Reset selection to A2 ---
$ws.Range("A2").Select() | Out-Null
